$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.412.74"
$ws.Range("E2").Value = "  +2.29%  "

$ws.Range("D3").Value = "2.237.40"
$ws.Range("E3").Value = "  +1.46%  "

$ws.Range("E4").Value = "  -0.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.79%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.583"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.09%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.563"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.40%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.34%  "

$ws.Range("E11").Value = "  -0.17%  "

$ws.Range("E12").Value = "  +3.28%  "

$ws.Range("E13").Value = "  +2.23%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.868"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.41%  "

$ws.Range("D15").Value = "2.576.58"
$ws.Range("E15").Value = "  +1.61%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.99%  "

$ws.Range("D17").Value = "2.236.81"
$ws.Range("E17").Value = "  +1.93%  "

$ws.Range("D18").Value = "43.330.03"
$ws.Range("E18").Value = "  +2.45%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.13"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.03%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.32%  "

$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "0.0₃0971"
$ws.Range("E21").Value = "  +2.63%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.23"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.80%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.66"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.37%  "

$ws.Range("E26").Value = "  +0.06%  "

$ws.Range("E27").Value = "  +3.41%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.20%  "

$ws.Range("E29").Value = "  +3.38%  "

$ws.Range("E30").Value = "  -1.82%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.73"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.96%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0881"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.45%  "

$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.30"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.39%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "156.52"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.48%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.72"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.57%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.23"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.65%  "

$ws.Range("E37").Value = "  +0.47%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.87"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.22%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.47%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.105"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.74%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.72"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.89%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0322"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.74%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "14.41"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +20.17%  "

$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("D45").Value = "1.805.22"
$ws.Range("E45").Value = "  +2.73%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.204"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.00%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "84.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.93"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.97%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.32%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.46%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "58.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.22%  "
